$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values (homework grades)
$ws.Range("C6").Value = 5
$ws.Range("D6").Value = 5

$ws.Range("F11").Value = 4

$ws.Range("F23").Value = 5
$ws.Range("G23").Value = 5

$ws.Range("C31").Value = 4
$ws.Range("D31").Value = 4
$ws.Range("E31").Value = 4
$ws.Range("F31").Value = 4

# Update view/selection state to match the saved workbook window
$ws.Range("G11").Select()
